# Progress in cohorts made:
# Extend the single BMI measurement row to also carry a second measurement
# (Weight) alongside the existing BMI one, by turning the relevant
# measures_* cells into comma-joined two-value strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AT2").Value = "LOINC:35925-4,LOINC:3141-9"
$ws.Range("AU2").Value = "BMI,Weight"
$ws.Range("AV2").Value = "24-9-2021,24-9-2021"
$ws.Range("BA2").Value = "NCIT:C49671,NCIT:C28252"
$ws.Range("BB2").Value = "Kilogram per Square Meter,Kilogram"
$ws.Range("BC2").Value = "26.63838307,85.6358"
